# Insert a new data row at row 412, shifting existing rows 412..494 down to 413..495.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(412).Insert()

# Populate the newly inserted row 412 with its data.
$ws.Range("A412").Value = 5
$ws.Range("B412").Value = "Macroferia Regional de Talca"
$ws.Range("C412").Value = "Maule"
$ws.Range("D412").Value = 44711
$ws.Range("E412").Value = 7
$ws.Range("F412").Value = "Fruta"
$ws.Range("G412").Value = 100109
$ws.Range("H412").Value = "Uva"
$ws.Range("I412").Value = 100109001
$ws.Range("J412").Value = "Uva"
$ws.Range("K412").Value = "Red Globe"
$ws.Range("L412").Value = "Primera"
$ws.Range("M412").Value = 150
$ws.Range("N412").Value = 10000
$ws.Range("O412").Value = 10000
$ws.Range("P412").Value = 10000
$ws.Range("Q412").Value = "$/bandeja 18 kilos"
$ws.Range("R412").Value = "Provincia de Limarí"
$ws.Range("S412").Value = 556
$ws.Range("T412").Value = 18
